$wb = $excel.ActiveWorkbook

# 1. hotel_info: insert a new "State" column right after "Hotel_Name" (before "City")
$hotelSheet = $wb.Worksheets.Item("hotel_info")
$hotelSheet.Columns.Item(3).Insert()
$hotelSheet.Range("C1").Value = "State"
$hotelSheet.Range("C2").Value = "Louisiana"

# 2. Reorder sheet tabs so "review_info" comes before "hotel_info"
$reviewSheet = $wb.Worksheets.Item("review_info")
$reviewSheet.Move($wb.Worksheets.Item(1))
